$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ligand = "Gnai2"
$receptor = "Ednra"

# Values per (sending cluster, target cluster) pair, in column order E..T
$valuesByPair = @{}
$valuesByPair["ECs|ECs"] = @(3,1,162.399297,487.197891,0.3910371682630009,0.3910371682630009,3,1,1.239942333333333,3.719827,0.02873264962564676,0.02873264962564676,201.365763253873,1812.291869284857,0.01123553394630589,0.01123553394630588)
$valuesByPair["ECs|FAPs"] = @(3,1,162.399297,487.197891,0.3910371682630009,0.3910371682630009,3,1,27.025359,81.076077,0.6262470038157576,0.6262470038157575,4388.899302772623,39500.0937249536,0.2448858550053026,0.2448858550053025)
$valuesByPair["ECs|M2"] = @(3,1,162.399297,487.197891,0.3910371682630009,0.3910371682630009,1,0.3333333333333333,0.050758,0.152274,0.001176193271648315,0.001176193271648314,8.243063517125998,74.18757165413399,0.0004599352862753516,0.0004599352862753515)
$valuesByPair["ECs|sCs"] = @(3,1,162.399297,487.197891,0.3910371682630009,0.3910371682630009,3,1,14.838413,44.51523899999999,0.3438441532869475,0.3438441532869475,2409.747839795661,21687.73055816094,0.1344558440251171,0.1344558440251171)
$valuesByPair["FAPs|ECs"] = @(3,1,65.41736466666667,196.252094,0.1575168212364948,0.1575168212364948,3,1,1.239942333333333,3.719827,0.02873264962564676,0.02873264962564676,81.11375978530423,730.0238380677379,0.004525875634733841,0.00452587563473384)
$valuesByPair["FAPs|FAPs"] = @(3,1,65.41736466666667,196.252094,0.1575168212364948,0.1575168212364948,3,1,27.025359,81.076077,0.6262470038157576,0.6262470038157575,1767.927764950582,15911.34988455524,0.09864443734993718,0.09864443734993716)
$valuesByPair["FAPs|M2"] = @(3,1,65.41736466666667,196.252094,0.1575168212364948,0.1575168212364948,1,0.3333333333333333,0.050758,0.152274,0.001176193271648315,0.001176193271648314,3.320454595750667,29.884091361756,0.0001852702253097956,0.0001852702253097955)
$valuesByPair["FAPs|sCs"] = @(3,1,65.41736466666667,196.252094,0.1575168212364948,0.1575168212364948,3,1,14.838413,44.51523899999999,0.3438441532869475,0.3438441532869475,970.6898742956072,8736.208868660466,0.05416123802651403,0.05416123802651403)
$valuesByPair["M2|ECs"] = @(3,1,126.3069433333333,378.92083,0.3041313008456065,0.3041313008456065,3,1,1.239942333333333,3.719827,0.02873264962564676,0.02873264962564676,156.6133260329345,1409.51993429641,0.008738498107388977,0.008738498107388977)
$valuesByPair["M2|FAPs"] = @(3,1,126.3069433333333,378.92083,0.3041313008456065,0.3041313008456065,3,1,27.025359,81.076077,0.6262470038157576,0.6262470038157575,3413.49048777599,30721.41438998391,0.1904613159211498,0.1904613159211498)
$valuesByPair["M2|M2"] = @(3,1,126.3069433333333,378.92083,0.3041313008456065,0.3041313008456065,1,0.3333333333333333,0.050758,0.152274,0.001176193271648315,0.001176193271648314,6.411087829713333,57.69979046742,0.0003577171897522517,0.0003577171897522517)
$valuesByPair["M2|sCs"] = @(3,1,126.3069433333333,378.92083,0.3041313008456065,0.3041313008456065,3,1,14.838413,44.51523899999999,0.3438441532869475,0.3438441532869475,1874.194589947596,16867.75130952837,0.1045737696273154,0.1045737696273155)
$valuesByPair["sCs|ECs"] = @(3,1,61.180387,183.541161,0.1473147096548978,0.1473147096548978,3,1,1.239942333333333,3.719827,0.02873264962564676,0.02873264962564676,75.86015181101632,682.7413662991469,0.004232741937218061,0.00423274193721806)
$valuesByPair["sCs|FAPs"] = @(3,1,61.180387,183.541161,0.1473147096548978,0.1473147096548978,3,1,27.025359,81.076077,0.6262470038157576,0.6262470038157575,1653.421922433933,14880.7973019054,0.092255395539368,0.092255395539368)
$valuesByPair["sCs|M2"] = @(3,1,61.180387,183.541161,0.1473147096548978,0.1473147096548978,1,0.3333333333333333,0.050758,0.152274,0.001176193271648315,0.001176193271648314,3.105394083346,27.948546750114,0.0001732705703109158,0.0001732705703109158)
$valuesByPair["sCs|sCs"] = @(3,1,61.180387,183.541161,0.1473147096548978,0.1473147096548978,3,1,14.838413,44.51523899999999,0.3438441532869475,0.3438441532869475,907.8198498058308,8170.378648252477,0.05065330160800083,0.05065330160800084)

$clusters = @("ECs","FAPs","M2","sCs")

$row = 2
foreach ($sc in $clusters) {
  foreach ($tc in $clusters) {
    $key = "$sc|$tc"
    $vals = $valuesByPair[$key]
    $ws.Cells.Item($row, 1).Value = $sc
    $ws.Cells.Item($row, 2).Value = $ligand
    $ws.Cells.Item($row, 3).Value = $receptor
    $ws.Cells.Item($row, 4).Value = $tc
    $col = 5
    foreach ($v in $vals) {
      $ws.Cells.Item($row, $col).Value = $v
      $col = $col + 1
    }
    $row = $row + 1
  }
}